$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.814.62"
$ws.Range("E2").Value = "  +3.57%  "

$ws.Range("D3").Value = "3.638.88"
$ws.Range("E3").Value = "  +2.20%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "572.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.23%  "

$ws.Range("D7").Value = "3.633.90"
$ws.Range("E7").Value = "  +2.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.627"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.39%  "

$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.685"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "62.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +17.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.153"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000292"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.96%  "

$ws.Range("D15").Value = "4.213.00"
$ws.Range("E15").Value = "  +2.11%  "

$ws.Range("D16").Value = "3.641.81"
$ws.Range("E16").Value = "  +2.30%  "

$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.93%  "

$ws.Range("D19").Value = "68.547.12"
$ws.Range("E19").Value = "  +3.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "407.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +19.98%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.89%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +15.81%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "669.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.116"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.424"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.35%  "

$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0812"
$ws.Range("E39").Value = "  +6.25%  "

$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.61%  "

$ws.Range("D42").Value = "3.222.32"
$ws.Range("E42").Value = "  +8.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.137"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.64%  "

$ws.Range("E44").Value = "  +11.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +25.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0423"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.133"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.81%  "
